$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "last charge end" timestamp applied to column D for rows 2-48
$ws.Range("D2:D48").Value = 45955.313148148147

# Refresh data rows 20-45 (row 19 is unchanged) with the latest pull; rows 46-48 are now empty
$ws.Range("A20").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B20").Value = "406号直流"
$ws.Range("C20").Value = 45952.398726851854
$ws.Range("A21").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B21").Value = "208号直流"
$ws.Range("C21").Value = 45953.419398148151
$ws.Range("A22").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B22").Value = "102号直流"
$ws.Range("C22").Value = 45953.530497685184
$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "505号直流"
$ws.Range("C23").Value = 45953.566435185188
$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "702号直流"
$ws.Range("C24").Value = 45953.571631944447
$ws.Range("A25").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B25").Value = "103号直流"
$ws.Range("C25").Value = 45953.589212962965
$ws.Range("A26").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B26").Value = "101号直流"
$ws.Range("C26").Value = 45954.028229166666
$ws.Range("A27").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value = "103号直流"
$ws.Range("C27").Value = 45954.036886574075
$ws.Range("A28").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B28").Value = "101号直流"
$ws.Range("C28").Value = 45954.071608796294
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "705号直流"
$ws.Range("C29").Value = 45954.245150462964
$ws.Range("A30").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B30").Value = "402号直流"
$ws.Range("C30").Value = 45954.274085648147
$ws.Range("A31").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B31").Value = "109号直流"
$ws.Range("C31").Value = 45954.323449074072
$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "B04号直流"
$ws.Range("C32").Value = 45954.460833333331
$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "B01号直流"
$ws.Range("C33").Value = 45954.481504629628
$ws.Range("A34").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B34").Value = "904号直流"
$ws.Range("C34").Value = 45954.539525462962
$ws.Range("A35").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value = "004A号直流"
$ws.Range("C35").Value = 45954.540092592593
$ws.Range("A36").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B36").Value = "108号直流"
$ws.Range("C36").Value = 45954.572337962964
$ws.Range("A37").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B37").Value = "106号直流"
$ws.Range("C37").Value = 45954.573819444442
$ws.Range("A38").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B38").Value = "105号直流"
$ws.Range("C38").Value = 45954.574155092596
$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "901号直流"
$ws.Range("C39").Value = 45954.582071759258
$ws.Range("A40").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B40").Value = "201号直流"
$ws.Range("C40").Value = 45954.583055555559
$ws.Range("A41").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B41").Value = "403号直流"
$ws.Range("C41").Value = 45954.607268518521
$ws.Range("A42").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value = "901号直流"
$ws.Range("C42").Value = 45954.725763888891
$ws.Range("A43").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B43").Value = "401号直流"
$ws.Range("C43").Value = 45954.741076388891
$ws.Range("A44").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B44").Value = "401号直流"
$ws.Range("C44").Value = 45954.747997685183
$ws.Range("A45").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B45").Value = "604号直流"
$ws.Range("C45").Value = 45954.77103009259

# Rows 46-48 no longer have data in this pull
$ws.Range("A46:C48").ClearContents()

# Selection moved to E22
$ws.Range("E22").Select() | Out-Null
